$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C (Fitness), rows 2 through 252, from 7293 to 7573
$ws.Range("C2:C252").Value = 7573
